$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F6").Value = 885
$ws1.Range("F9").Value = 1294
$ws1.Range("F11").Value = 3183
$ws1.Range("F14").Value = 1205
$ws1.Range("F18").Value = 1818
$ws1.Range("F20").Value = 555376
$ws1.Range("F22").Value = 270
$ws1.Range("F23").Value = 622
$ws1.Range("F26").Value = 1710
$ws1.Range("F29").Value = 491
$ws1.Range("F30").Value = 1335
$ws1.Range("F38").Value = 1080
$ws1.Range("F41").Value = 2356
$ws1.Range("F44").Value = 2875
$ws1.Range("F46").Value = 855

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F14").Value = 126771
$ws2.Range("F15").Value = 126771
$ws2.Range("F21").Value = 213
$ws2.Range("F25").Value = 7
$ws2.Range("F28").Value = 74
$ws2.Range("F29").Value = 224
$ws2.Range("G31").Value = 380
$ws2.Range("F35").Value = 187
$ws2.Range("F40").Value = 96

$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F6").Value = 4864
$ws3.Range("F9").Value = 717
$ws3.Range("F10").Value = 1002
$ws3.Range("F13").Value = 1449
$ws3.Range("F15").Value = 1426

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 717
$ws4.Range("F5").Value = 1002
$ws4.Range("F8").Value = 1449
$ws4.Range("F12").Value = 1294
$ws4.Range("F13").Value = 1426
$ws4.Range("F14").Value = 3183
$ws4.Range("F18").Value = 1205
$ws4.Range("F21").Value = 1818
$ws4.Range("F24").Value = 555377
$ws4.Range("F27").Value = 126771
$ws4.Range("F29").Value = 622
$ws4.Range("F34").Value = 491
$ws4.Range("F35").Value = 1335
$ws4.Range("G43").Value = 380
$ws4.Range("F44").Value = 1080
$ws4.Range("F45").Value = 187
$ws4.Range("F47").Value = 2356
